# Applies the 'Busca de produtos por nome e categorias' edit described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch column A so the sheet's used-range/dimension grows to start at A1
# (matches the saved file, whose dimension goes from B1:G106 to A1:G117).
$ws.Range('A1').Value = 'x'
$ws.Range('A1').Value = $null

# --- 1) Existing rows 102-106: insert a new E-column title, push the old
#        'nome aula' text for that row down into the F column. ---------------
$ws.Range('F102').Value = '0:44 - discussão sobre Caso de Uso - descreve o cenário de utilização da aplicação...'
$ws.Range('E102').Value = 'Apresentando o caso de uso'

$ws.Range('F103').Value = '0:34 - JPQL é a linguagem de consulta da JPA'
$ws.Range('E103').Value = 'Nivelamento sobre SQL e JPQL'

$ws.Range('F104').Value = '1:27 - exemplo/comparativo entre uma consulta feita com SQL e uma consulta feita com JPQL'
$ws.Range('E104').Value = 'Nivelamento sobre SQL e JPQL'

$ws.Range('F105').Value = '
9:31 - sintaxe INNER JOIN em SQl e JPQL'
$ws.Range('E105').Value = 'Nivelamento sobre SQL e JPQL'

$ws.Range('F106').Value = '
12:48 - instrução SELECT DISTINCT - faz uma consulta no banco de dados e retorna objetos sem repetição'
$ws.Range('E106').Value = 'Nivelamento sobre SQL e JPQL'

# --- 2) Grow the table ('Tabela1') by 11 rows (107-117) ----------------------
$lo = $ws.ListObjects.Item('Tabela1')
for ($i = 0; $i -lt 11; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# --- 3) Fill in the 11 new rows ------------------------------------------------
$ws.Range('B107').Value = 3
$ws.Range('C107').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D107').Value = 50
$ws.Range('E107').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F107').Value = '4:20 - criação de um novo service, no caso ProdutoService'
$ws.Range('G107').Value = '
'

$ws.Range('B108').Value = 3
$ws.Range('C108').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D108').Value = 50
$ws.Range('E108').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F108').Value = '5:22 - criação do metodo search com paginação da classe ProdutoService'

$ws.Range('B109').Value = 3
$ws.Range('C109').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D109').Value = 50
$ws.Range('E109').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F109').Value = '9:26 - Spring Data - na documentação mostra os "Query Methods"... são metodos que utilizam da nomeação padrao de nomes de metodos para gerar automaticamente a lógica de consulta ... util pois nao precisa implementar alguns metodos'

$ws.Range('B110').Value = 3
$ws.Range('C110').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D110').Value = 50
$ws.Range('E110').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F110').Value = '10:30 - anotação @Query("INSIRA SUA JPQL") do framework Spring Data- utilizada na interface de repositórios. o framework faz um pré-processamento e cria automaticamente o metodo sem precisar a criar uma nova classe'

$ws.Range('B111').Value = 3
$ws.Range('C111').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D111').Value = 50
$ws.Range('E111').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F111').Value = '10:41 - anotação @Param - responsável por pegar o valor do parâmetro passado no método (search) com a anotação @Query e passar para o parâmetro da consulta JPQL... neste caso, nomeamos o parâmetro da JPQL como "%:nome%"'

$ws.Range('B112').Value = 3
$ws.Range('C112').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D112').Value = 50
$ws.Range('E112').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F112').Value = '12:42 - criação do end point ProdutoResource'

$ws.Range('B113').Value = 3
$ws.Range('C113').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D113').Value = 50
$ws.Range('E113').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F113').Value = '14:11 - criação da classe ProdutoDTO'

$ws.Range('B114').Value = 3
$ws.Range('C114').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D114').Value = 50
$ws.Range('E114').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F114').Value = '15:20 - o metodo GET não aceita enviar parametros no corpo da requisição .. como o POST ... somente aceita como parametros na URL'

$ws.Range('B115').Value = 3
$ws.Range('C115').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D115').Value = 50
$ws.Range('E115').Value = 'Busca de pedidos por nome e categorias - PARTE 1'
$ws.Range('F115').Value = '18:11 - criação de classe auxiliar URL'

$ws.Range('B116').Value = 3
$ws.Range('C116').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D116').Value = 51
$ws.Range('E116').Value = 'Busca de pedidos por nome e categorias - PARTE 2'
$ws.Range('F116').Value = '5:37 - F A N T Á S T I C O - dica de uso do padrao de nomes do framework descrito na documentação, para que o próprio framework se encarregue de gerar as consultas JPQL invés de criarmos as consultas/sintaxes na unha'
$r2 = $ws.Range('F116').Characters(8, 19)
$r2.Font.Bold = $true
$r2.Font.Size = 14
$r2.Font.Color = 255
$r2.Font.Name = 'Consolas'
$r3 = $ws.Range('F116').Characters(27, 190)
$r3.Font.Size = 11
$r3.Font.Name = 'Consolas'
$ws.Range('G116').Value = '

'

$ws.Range('B117').Value = 3
$ws.Range('C117').Value = 'Operações de CRUD e Casos de Uso'
$ws.Range('D117').Value = 51
$ws.Range('E117').Value = 'Busca de pedidos por nome e categorias - PARTE 2'
$ws.Range('F117').Value = '1:52 - criação de metodo de encode - que formata a URL eliminando caracteres invalidos, como por exemplo, espaços...'

# --- 4) Cosmetic: zoom + final selection, matching the saved view -------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range('F117').Select()
